# Update Bibi faturamento diario lojas data:
#  - Swap the store names "Bibi Cell Vieiralves" and "Bibi Cell Ponta Negra"
#    (row 3 and row 4 in column A)
#  - Replace the daily values for rows 2-6 (stores + total) with a single
#    value in "day 1" (column B), zero out the other day columns (C..AE),
#    and update the row total (column AG) to match the new day-1 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap store names in column A (rows 3 and 4) ---
$ws.Range("A3").Value = "Bibi Cell Ponta Negra"
$ws.Range("A4").Value = "Bibi Cell Vieiralves"

# --- New day-1 (column B) values per row, and matching new totals (column AG) ---
$newValues = @{
    2 = 6805.15
    3 = 4535.01
    4 = 3638
    5 = 2251
    6 = 17229.16
}

foreach ($row in $newValues.Keys) {
    $value = $newValues[$row]

    # Zero out columns C..AE (day 2 .. day 30)
    $ws.Range("C" + $row + ":AE" + $row).Value = 0

    # Set new day-1 value (column B)
    $ws.Range("B" + $row).Value = $value

    # Update the row total (column AG)
    $ws.Range("AG" + $row).Value = $value
}
